$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.843.08'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.931.00'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'352.59"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = "'111.87"
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").Value = "'39.33"
$ws.Range("E10").Value = '  -2.17%  '
$ws.Range("E11").Value = '  +4.67%  '
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = "'20.04"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = "'7.78"
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").Value = '3.393.27'
$ws.Range("E15").Value = '  +3.30%  '
$ws.Range("D16").Value = '2.934.86'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '51.913.19'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("E19").Value = '  -4.62%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = "'14.21"
$ws.Range("E21").Value = '  +6.47%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").Value = "'71.24"
$ws.Range("D24").Value = "'268.53"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E26").Value = '  +11.54%  '
$ws.Range("D27").Value = "'26.91"
$ws.Range("E27").Value = '  +2.42%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = '  +15.43%  '
$ws.Range("E30").Value = '  +15.67%  '
$ws.Range("D31").Value = "'10.56"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = "'6.29"
$ws.Range("E32").Value = '  +11.38%  '
$ws.Range("D33").Value = "'2.27"
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("D34").Value = "'37.16"
$ws.Range("E34").Value = '  -4.62%  '
$ws.Range("D35").Value = "'52.96"
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  +3.46%  '
$ws.Range("D39").Value = "'18.76"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = '  +6.17%  '
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("D43").Value = "'23.17"
$ws.Range("E43").Value = '  +3.96%  '
$ws.Range("E44").Value = '  -1.07%  '
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").Value = "'3.51"
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").Value = '2.174.67'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = "'111.63"
$ws.Range("E48").Value = '  -8.82%  '
$ws.Range("E49").Value = '  +2.59%  '
$ws.Range("D50").Value = "'0.0349"
$ws.Range("E50").Value = '  +10.73%  '
$ws.Range("D51").Value = "'0.943"
$ws.Range("E51").Value = '  -1.38%  '
